$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $ws.Range("E$r").Value = $ws.Range("E$r").Value() - 1
    $ws.Range("G$r").Value = $ws.Range("G$r").Value() + 1
    $ws.Range("AQ$r").Value = $ws.Range("AQ$r").Value() - 1
    $ws.Range("AS$r").Value = $ws.Range("AS$r").Value() + 1
}
